# Populate the four worksheets (Products, Customers, Employees, Invoices)
# with header rows and two sample data rows each, matching the target
# OOXML diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Products
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A1").Value = "id"
$ws1.Range("B1").Value = "name"
$ws1.Range("C1").Value = "sku"
$ws1.Range("D1").Value = "category"
$ws1.Range("E1").Value = "price"
$ws1.Range("F1").Value = "cost_price"
$ws1.Range("G1").Value = "stock_quantity"
$ws1.Range("H1").Value = "unit"
$ws1.Range("I1").Value = "hsn_code"
$ws1.Range("J1").Value = "gst_rate"
$ws1.Range("K1").Value = "is_active"

$ws1.Range("A2").Value = "594e2ff5-1dac-4141-bfbf-0888c697862a"
$ws1.Range("B2").Value = "Alpha Widget"
$ws1.Range("C2").Value = "AW-100"
$ws1.Range("D2").Value = "Widgets"
$ws1.Range("E2").Value = 99.99
$ws1.Range("F2").Value = 70
$ws1.Range("G2").Value = 20
$ws1.Range("H2").Value = "piece"
$ws1.Range("I2").NumberFormat = "@"
$ws1.Range("I2").Value = "1001"
$ws1.Range("J2").Value = 18
$ws1.Range("K2").Value = $true

$ws1.Range("A3").Value = "d6647f0d-be66-49c0-9947-304d0d13eaad"
$ws1.Range("B3").Value = "Beta Gadget"
$ws1.Range("C3").Value = "BG-200"
$ws1.Range("D3").Value = "Gadgets"
$ws1.Range("E3").Value = 149.5
$ws1.Range("F3").Value = 120
$ws1.Range("G3").Value = 10
$ws1.Range("H3").Value = "piece"
$ws1.Range("I3").NumberFormat = "@"
$ws1.Range("I3").Value = "2002"
$ws1.Range("J3").Value = 12
$ws1.Range("K3").Value = $true

# ---------------------------------------------------------------------
# Sheet 2: Customers
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1").Value = "id"
$ws2.Range("B1").Value = "name"
$ws2.Range("C1").Value = "email"
$ws2.Range("D1").Value = "phone"
$ws2.Range("E1").Value = "gstin"
$ws2.Range("F1").Value = "address"

$ws2.Range("A2").Value = "b4059366-584f-45d1-b842-2e2054e9b7f1"
$ws2.Range("B2").Value = "John Doe"
$ws2.Range("C2").Value = "john@example.com"
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "1234567890"
$ws2.Range("E2").Formula = '=""'
$ws2.Range("F2").Value = "123 Road City"

$ws2.Range("A3").Value = "38269074-e958-4faf-bd21-474a68066b72"
$ws2.Range("B3").Value = "Acme Corp"
$ws2.Range("C3").Value = "acme@corp.com"
$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = "5550002222"
$ws2.Range("E3").Value = "27AACCA1234F1Z2"
$ws2.Range("F3").Value = "456 Industrial Area"

# ---------------------------------------------------------------------
# Sheet 3: Employees
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A1").Value = "id"
$ws3.Range("B1").Value = "name"
$ws3.Range("C1").Value = "title"
$ws3.Range("D1").Value = "salary"

$ws3.Range("A2").Value = "92653c1f-0553-4d9a-a9d9-6d5942787f97"
$ws3.Range("B2").Value = "Jane Smith"
$ws3.Range("C2").Value = "Manager"
$ws3.Range("D2").Value = 60000

$ws3.Range("A3").Value = "8e8d9d86-43b8-421b-87ba-1a743b5d61ca"
$ws3.Range("B3").Value = "Bob Miller"
$ws3.Range("C3").Value = "Staff"
$ws3.Range("D3").Value = 30000

# ---------------------------------------------------------------------
# Sheet 4: Invoices
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("A1").Value = "id"
$ws4.Range("B1").Value = "customer_id"
$ws4.Range("C1").Value = "total"
$ws4.Range("D1").Value = "created_at"
$ws4.Range("E1").Value = "items"

$ws4.Range("A2").Value = "c6443c74-26d6-4a65-b20f-20ad1e015232"
$ws4.Range("B2").Value = "b4059366-584f-45d1-b842-2e2054e9b7f1"
$ws4.Range("C2").Value = 500
$ws4.Range("D2").Value = "2025-10-30T14:54:07.153Z"
$ws4.Range("E2").Value = '[{"product_id":"594e2ff5-1dac-4141-bfbf-0888c697862a","qty":2,"price":99.99}]'

$ws4.Range("A3").Value = "5ecf3b24-9c3a-4e8c-9b99-aede363b2077"
$ws4.Range("B3").Value = "38269074-e958-4faf-bd21-474a68066b72"
$ws4.Range("C3").Value = 149.5
$ws4.Range("D3").Value = "2025-10-30T14:54:07.154Z"
$ws4.Range("E3").Value = '[{"product_id":"d6647f0d-be66-49c0-9947-304d0d13eaad","qty":1,"price":149.5}]'
